$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 9: "31.05.2024" / "Rohkem teksti. ..." ---
# (written first so these strings land at shared-string indices 11/12,
#  matching the author's original edit order)
$ws.Range("A9").Value = "31.05.2024"
$ws.Range("B9").Value = "Rohkem teksti.                            Lisatud emotsioone mõne tegelase jaoks.                             Tähemärkide lisamine koodile"

# --- Row 7: extend the existing note with the new sentence about emotions ---
$ws.Range("B7").Value = "Jätkasime dialooge. Taustapiltide tegemine.          Joonistasin tegelastele emotsioone."

# --- New row 8: "28.05.2024" / "Joonistasin tegelastele emotsioone." ---
$ws.Range("A8").Value = "28.05.2024"
$ws.Range("B8").Value = "Joonistasin tegelastele emotsioone."

# Match the wrap-text formatting used throughout column B (and now A9)
$ws.Range("B8").WrapText = $true
$ws.Range("A9").WrapText = $true
$ws.Range("B9").WrapText = $true

# Row heights to fit the new/changed text
$ws.Rows(7).RowHeight = 60
$ws.Rows(8).RowHeight = 30
$ws.Rows(9).RowHeight = 75

# Update selection / scroll position
$ws.Range("F11").Select()
